# Testati con gpt tutti i progetti
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Totals header row (row 19): new column headers for the extra breakdown columns ---
# (written first so new shared-string order matches the source workbook)
$ws.Range("E19").Value = "Generazioni mancanti"
$ws.Range("F19").Value = "Generazioni non necessarie"
$ws.Range("D19").Value = "Totale fallimenti"
$ws.Range("G19").Value = "Generazioni necessarie ma errate"

# --- Title cell (A1): project name changed to "Progetto-Alfredo" ---
$ws.Range("A1").Value = "Applicazione: Progetto-Alfredo - Confronto Robustezza Locatori"

# --- LLM block (rows 4-9): fill in the real results that used to be "N\D" placeholders ---
# hook
$ws.Range("C4").Value = 44
$ws.Range("D4").Value = 42
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Formula = "=(E4/C4)*100"
# absolute
$ws.Range("C5").Value = 44
$ws.Range("D5").Value = 33
$ws.Range("E5").Value = 11
$ws.Range("F5").Value = 1
$ws.Range("G5").Formula = "=(E5/C5)*100"
# relative
$ws.Range("C6").Value = 44
$ws.Range("D6").Value = 36
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 1
$ws.Range("G6").Formula = "=(E6/C6)*100"
# robula
$ws.Range("C7").Value = 44
$ws.Range("D7").Value = 40
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Formula = "=(E7/C7)*100"
# selenium
$ws.Range("C8").Value = 44
$ws.Range("D8").Value = 36
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 1
$ws.Range("G8").Formula = "=(E8/C8)*100"
# katalon
$ws.Range("C9").Value = 44
$ws.Range("D9").Value = 40
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Formula = "=(E9/C9)*100"

# Formatting for the newly filled LLM cells: "Totale Test"/"Fallimenti"/"Obsolescenza"
# columns centered both horizontally & vertically, tasso column kept as a centered 2-decimal %.
$ws.Range("C4:C9").HorizontalAlignment = -4108
$ws.Range("C4:C9").VerticalAlignment = -4108
$ws.Range("E4:E9").HorizontalAlignment = -4108
$ws.Range("E4:E9").VerticalAlignment = -4108
$ws.Range("F4:F9").HorizontalAlignment = -4108
$ws.Range("F4:F9").VerticalAlignment = -4108

$ws.Range("G4:G9").NumberFormat = "0.00"
$ws.Range("G4:G9").HorizontalAlignment = -4108

# --- LLM totals row (row 20) ---
$ws.Range("B20").Formula = "=SUM(E4:E9)"
$ws.Range("C20").Formula = "=SUM(F4:F9)"
$ws.Range("D20").Formula = "=SUM(B20,C20)"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 0

# --- Analitica totals row (row 21): fill the extra breakdown columns ---
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0
$ws.Range("F21").HorizontalAlignment = -4108
$ws.Range("G21").Value = 0

# --- Misc view state tweaks ---
$ws.Columns.Item(7).AutoFit()
$ws.Range("D7").Select()
